# feat: add 2022-Q4 data
#
# The workbook currently has two sheets: "总计" (summary) and "2021-Q4"
# (fund holdings detail for Q4 2021). This change inserts a new
# "2022-Q4" detail sheet (positioned right after "总计", before
# "2021-Q4"), updates the summary sheet so its single data row now
# refers to 2022-Q4, and appends a new summary row for the pre-existing
# 2021-Q4 data.

$wb = $excel.ActiveWorkbook

# --- 1. Create the new "2022-Q4" detail sheet -----------------------
# Duplicate the existing "2021-Q4" sheet (same columns/header/styling)
# so the new sheet inherits identical formatting, then drop it right
# after "总计" and rename it.
$total = $wb.Worksheets.Item("总计")
$oldQ4 = $wb.Worksheets.Item("2021-Q4")
$oldQ4.Copy($null, $total)

$newQ4 = $wb.Worksheets.Item("2021-Q4 (2)")
$newQ4.Name = "2022-Q4"

# Update the figures that differ for the 2022-Q4 snapshot. Columns D:G
# are stored as text (e.g. "82.94", "0.0037") in the source sheet, so
# force text format before assigning, otherwise Excel would coerce
# these numeric-looking strings into real numbers.
$newQ4.Range("D2:G3").NumberFormat = "@"

$newQ4.Range("C2").Value = "长信美国标准普尔100等权重指数增强（QDII）人民币"
$newQ4.Range("D2").Value = "0.44"
$newQ4.Range("E2").Value = "82.94"
$newQ4.Range("F2").Value = "0.85"
$newQ4.Range("G2").Value = "0.0037"

$newQ4.Range("C3").Value = "长信美国标准普尔100等权重指数增强（QDII）美元"
$newQ4.Range("D3").Value = "0.44"
$newQ4.Range("E3").Value = "82.94"
$newQ4.Range("F3").Value = "0.85"
$newQ4.Range("G3").Value = "0.0037"

# --- 2. Update the "总计" summary sheet ------------------------------
# Row 2 used to describe 2021-Q4; it now describes 2022-Q4 (counts are
# unchanged). Row 3 is a new row that preserves the original 2021-Q4
# summary figures, copying A2's formatting (bold/bordered index style)
# onto the new A3 cell.
$total.Range("A2").Copy($total.Range("A3"))
$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 2
$total.Range("D3").Value = 0.01

$total.Range("B2").Value = "2022-Q4"
